# Se cargo la base de datos de informacion de AUTOGENERACION
# Adds a new "total_2050" column (H) to the TimePeriods sheet and
# makes TimePeriods the active sheet/selection instead of Constants.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TimePeriods")

# --- New header cell H35: "total_2050" ------------------------------------
# Give it the same value as the neighbouring header cells and copy the
# formatting (fill/border/font) from G35 so it matches the existing header
# row style (s="66").
$ws.Range("H35").Value = "total_2050"
$ws.Range("G35").Copy()
$ws.Range("H35").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- New data column H37:H68: sequential years 2019-2050 -------------------
for ($i = 0; $i -lt 32; $i++) {
    $row = 37 + $i
    $ws.Cells.Item($row, 8).Value = 2019 + $i
}

# --- Make TimePeriods the active sheet with the new range selected ---------
$ws.Activate() | Out-Null
$ws.Range("H58:H68").Select() | Out-Null
